# Workbook has sheets: Setup, SetupTests, LaunchPage, LaunchPageScript, Login, LoginScript
$wb = $excel.ActiveWorkbook

# Update the phone-number test data on the Login sheet: A2 changes from a
# numeric value to a text string "7792016895".
$loginSheet = $wb.Worksheets.Item("Login")
$loginSheet.Range("A2").Value = "7792016895"

# Select cell A3 on the Login sheet and make that sheet the active tab.
$loginSheet.Activate()
$loginSheet.Range("A3").Select()
